# Updates TC08_Canine_Filter_Study-TCL01 'startup' worksheet for the
# TCL01 Study 1-101 test cases: each Cypher query gains an
# "and demo.breed in ['Beagle']" filter (plus order/limit clauses), the
# FilesTab query is rewritten, and the now-unused 'cartQuery' column is
# removed so the sheet holds TabName/query/StatQuery/dbExcel/WebExcel (A:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = 'TabName'
$ws.Range("B1").Value = 'query'
$ws.Range("C1").Value = 'StatQuery'
$ws.Range("D1").Value = 'dbExcel'
$ws.Range("E1").Value = 'WebExcel'
$ws.Range("A2").Value = 'CasesTab'
$ws.Range("B2").Value = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
 WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in ['Beagle']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc 
limit 100
'@
$ws.Range("C2").Value = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
 WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in ['Beagle']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$ws.Range("D2").Value = 'TC08_Canine_Filter_Study-TCL01_Neo4jData.xlsx'
$ws.Range("E2").Value = 'TC08_Canine_Filter_Study-TCL01_WebData.xlsx'
$ws.Range("A3").Value = 'SamplesTab'
$ws.Range("B3").Value = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
  WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in ['Beagle']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc 
limit 100
'@
$ws.Range("C3").Value = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
 WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in ['Beagle']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$ws.Range("D3").Value = 'TC08_Canine_Filter_Study-TCL01_Neo4jData.xlsx'
$ws.Range("E3").Value = 'TC08_Canine_Filter_Study-TCL01_WebData.xlsx'
$ws.Range("A4").Value = 'FilesTab'
$ws.Range("B4").Value = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in['Beagle']  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
order by f.filename asc
limit 100
'@
$ws.Range("C4").Value = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
 WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in ['Beagle']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$ws.Range("D4").Value = 'TC08_Canine_Filter_Study-TCL01_Neo4jData.xlsx'
$ws.Range("E4").Value = 'TC08_Canine_Filter_Study-TCL01_WebData.xlsx'
$ws.Range("A5").Value = 'StudyFilesTab'
$ws.Range("B5").Value = @'
  MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
 WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in ['Beagle']
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit

        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
order by f.filename asc
limit 100
'@
$ws.Range("C5").Value = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
 WHERE s.clinical_study_designation IN ['TCL01'] and demo.breed in ['Beagle']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@
$ws.Range("D5").Value = 'TC08_Canine_Filter_Study-TCL01_Neo4jData.xlsx'
$ws.Range("E5").Value = 'TC08_Canine_Filter_Study-TCL01_WebData.xlsx'

# Drop the now-empty column F (old WebData filename slot; that value now
# lives in column E since the 'cartQuery' column was collapsed away).
$ws.Range("F1:F5").Delete()

# Column widths (characters), matching the resized layout.
$ws.Columns.Item(1).ColumnWidth = 13.42
$ws.Columns.Item(2).ColumnWidth = 57.27
$ws.Columns.Item(3).ColumnWidth = 50.42
$ws.Columns.Item(4).ColumnWidth = 20.58
$ws.Columns.Item(5).ColumnWidth = 19.75
$ws.Columns.Item(6).ColumnWidth = 30.75

# Row heights (points), matching the wrapped-text layout for the new content.
$ws.Rows.Item(1).RowHeight = 29
$ws.Rows.Item(2).RowHeight = 377
$ws.Rows.Item(3).RowHeight = 304.5
$ws.Rows.Item(4).RowHeight = 275.5
$ws.Rows.Item(5).RowHeight = 409.5

# Match the saved selection/active cell.
$ws.Range("E5").Select()
